$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the 2020 column (Q) header
$ws.Range("Q3").Value = 2020

# Fill column Q with "-" for data rows 4-13 (mirrors the existing "-" placeholders)
$dash = "-"
$ws.Range("Q4").Value = $dash
$ws.Range("Q5").Value = $dash
$ws.Range("Q6").Value = $dash
$ws.Range("Q7").Value = $dash
$ws.Range("Q8").Value = $dash
$ws.Range("Q9").Value = $dash
$ws.Range("Q10").Value = $dash
$ws.Range("Q11").Value = $dash
$ws.Range("Q12").Value = $dash
$ws.Range("Q13").Value = $dash

# Copy formatting from column P so the new column matches existing styling
$ws.Range("P3:P13").Copy()
$ws.Range("Q3:Q13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the values (paste-formats only should have kept them, but make sure)
$ws.Range("Q3").Value = 2020
$ws.Range("Q4").Value = $dash
$ws.Range("Q5").Value = $dash
$ws.Range("Q6").Value = $dash
$ws.Range("Q7").Value = $dash
$ws.Range("Q8").Value = $dash
$ws.Range("Q9").Value = $dash
$ws.Range("Q10").Value = $dash
$ws.Range("Q11").Value = $dash
$ws.Range("Q12").Value = $dash
$ws.Range("Q13").Value = $dash

# Move/activate selection to mirror the recorded cursor position in the source file
$ws.Range("P17").Select()
